$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.05459999999999
$ws.Range("B9").Value = 8.453300000000002
$ws.Range("C11").Value = -13.5981
$ws.Range("B18").Value = 4.548500000000004
$ws.Range("B20").Value = 5.762599999999997
$ws.Range("D21").Value = -7.525400000000001
